$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append 45 new data rows (rows 102-146), following the same repeating
# pattern as the existing data: column A cycles 10002..10010, column B
# increments sequentially from 3000121, column C = "eng", column D = TRUE,
# column E = "superadmin", column F = "now()".
$startRow = 102
$startB = 3000121
$count = 45

for ($i = 0; $i -lt $count; $i++) {
    $row = $startRow + $i
    $aVal = 10002 + ($i % 9)
    $bVal = $startB + $i

    $ws.Cells.Item($row, 1).Value = $aVal
    $ws.Cells.Item($row, 2).Value = $bVal
    $ws.Cells.Item($row, 3).Value = "eng"
    $ws.Cells.Item($row, 4).Value = $true
    $ws.Cells.Item($row, 5).Value = "superadmin"
    $ws.Cells.Item($row, 6).Value = "now()"
}

# Update the view: scroll so row 128 is the top-left visible row, and
# select the range covering the newly added rows.
$ws.Application.ActiveWindow.ScrollRow = 128
$ws.Range("A102:F146").Select() | Out-Null

# Configure the page setup for printing (portrait orientation).
$ws.PageSetup.Orientation = 1
